$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws3 = $wb.Worksheets.Item("CONVERTION")

# --- Insert a new row at sheet row 77 (shifts old rows 77..123 down to 78..124) ---
$ws.Rows.Item(77).Insert()

# Fix up the formatting of the freshly-inserted row 77 so it matches the
# surrounding table rows (Insert() drops in a blank/default style otherwise).
$ws.Range("A76:K76").Copy()
$ws.Range("A77:K77").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Restore the calculated-column formula in G77 (new rows don't inherit it automatically)
$ws.Range("G77").Formula = "=IF(ISBLANK(Table13[[#This Row],[EARNED]]),"""",Table13[[#This Row],[EARNED]])"

# Grow the Excel Table (Table13) so it covers the new last row too
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A8:K124"))

# Restore the calculated-column formula on the new last (totals-style) row 124 as well
$ws.Range("G124").Formula = "=IF(ISBLANK(Table13[[#This Row],[EARNED]]),"""",Table13[[#This Row],[EARNED]])"

# --- New row 77 content: "UT(3-1-13)" undertime entry ---
$ws.Range("B77").Value2 = "UT(3-1-13)"
$ws.Range("D77").Value2 = 3.152

# --- Row 73 gets a new "UT(0-0-3)" undertime entry ---
$ws.Range("B73").Value2 = "UT(0-0-3)"
$ws.Range("D73").Value2 = 0.006

# --- CONVERTION sheet: move the "3" from E3 to F3 ---
$ws3.Range("E3").ClearContents()
$ws3.Range("F3").Value2 = 3

# --- Update the view selection to match ---
$ws.Activate()
$ws.Range("F80").Select()
